$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# Row 3: the Admission No (C3) was originally stored as text "14" and gets
# normalized to a real number as part of this sync.
$ws.Range("C3").Value = 14

# Row 4: a brand-new form submission synced in.
$ws.Range("A4").Value = "2026-02-08 03:56:36"
$ws.Range("B4").Value = "SARAH MUSA BALAMI"

# Admission No for the new row stays textual ("38"), matching the source
# form data, so force text storage and then reset the style back to the
# sheet default (no explicit cell style) to match the rest of the column.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "38"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = 10
